# Add a new "2021" column (R) to the right of the existing "2020" column (Q),
# mirroring the Q column's per-row number formatting/styles, then fill in the
# 2021 data values published alongside the existing 2007-2020 series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format / font / borders / etc.) from column Q
# (rows 3-33, the data rows) onto the new column R so each row picks up the
# same style the corresponding Q cell already uses.
$ws.Range("Q3:Q33").Copy()
$ws.Range("R3:R33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3 holds the year headers; rows 4-33 hold the yearly data values.
$rValues = @{
    3  = 2021
    4  = 1.7931687443515183
    5  = 1.0977143806517458
    6  = 2.4989281705678046
    7  = 2.3489023398681002
    8  = 1.8410239038543676
    9  = 2.8382683724659588
    10 = 1.2584206034913306
    11 = 0.79202525610136665
    12 = 1.7183687369364922
    13 = 1.7860084101151579
    14 = 1.5807090270340762
    15 = 1.9930959157478496
    16 = 1.0231016349164126
    17 = 0
    18 = 2.0091214112068791
    19 = 2.2092990108041848
    20 = 0.86496336159360854
    21 = 3.5236628052020538
    22 = 1.4678252700798498
    23 = 0.74155920237892192
    24 = 2.1792664589099311
    25 = 1.5302890103825006
    26 = 0.80351618683358383
    27 = 2.280288974802807
    28 = 2.3014726663297309
    29 = 1.7358308467556451
    30 = 2.9402079315049163
    31 = 1.2198989923634325
    32 = 1.1878318505232399
    33 = 1.2537455648750642
}

foreach ($row in $rValues.Keys) {
    $ws.Cells.Item($row, 18).Value = $rValues[$row]
}

# Match the author's final selection (cell S14, just past the new data).
$ws.Range("S14").Select()
